# LibraryReleases.xlsx — re-order the Type categories.
# The "Type" column now carries an explicit numeric prefix
# (1. operating system, 2. web server, 3. database, 4. framework,
#  5. cms, 6. library) and the table rows are re-sorted into that
# new priority order instead of the old alphabetical order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Re-write the B:D data rows (2-16) in the new sort order ----------
# Columns: B = Type, C = Product, D = url - version/release (display text;
# the real hyperlink target is set further down since a couple of them
# differ from the displayed text, e.g. the jQuery releases link).

$rows = @(
    @{ Row = 2;  Type = "1. operating system"; Product = "MS Windows Server"; Url = "https://learn.microsoft.com/en-us/windows-server/get-started/windows-server-release-info" },
    @{ Row = 3;  Type = "1. operating system"; Product = "Ubuntu";             Url = "https://wiki.ubuntu.com/Releases" },
    @{ Row = 4;  Type = "2. web server";       Product = "Apache";             Url = "https://httpd.apache.org/download.cgi" },
    @{ Row = 5;  Type = "2. web server";       Product = "IIS";                Url = "https://learn.microsoft.com/en-us/lifecycle/products/internet-information-services-iis" },
    @{ Row = 6;  Type = "2. web server";       Product = "Nginx";              Url = "https://nginx.org/en/download.html" },
    @{ Row = 7;  Type = "3. database";         Product = "MariaDB";            Url = "https://mariadb.com/kb/en/mariadb-server-release-dates/" },
    @{ Row = 8;  Type = "3. database";         Product = "MS SQL";             Url = "https://learn.microsoft.com/en-us/troubleshoot/sql/releases/download-and-install-latest-updates" },
    @{ Row = 9;  Type = "3. database";         Product = "MySQL";              Url = "https://dev.mysql.com/doc/refman/8.4/en/mysql-releases.html" },
    @{ Row = 10; Type = "4. framework";        Product = ".Net";               Url = "https://dotnet.microsoft.com/en-us/platform/support/policy/dotnet-core/" },
    @{ Row = 11; Type = "4. framework";        Product = "PHP";                Url = "https://www.php.net/supported-versions.php" },
    @{ Row = 12; Type = "5. cms";              Product = "Moodle";             Url = "https://moodledev.io/general/releases" },
    @{ Row = 13; Type = "5. cms";              Product = "Wordpress";          Url = "https://wordpress.org/download/releases/" },
    @{ Row = 14; Type = "6. library";          Product = "Bootstrap";          Url = "https://getbootstrap.com/docs/versions" },
    @{ Row = 15; Type = "6. library";          Product = "jQuery";             Url = "https://releases.jquery.com" },
    @{ Row = 16; Type = "6. library";          Product = "PHPMyAdmin";         Url = "https://www.phpmyadmin.net/downloads/" }
)

# The jQuery "releases.jquery.com" hyperlink target has a trailing slash
# even though the displayed text in D15 does not.
$hyperlinkTargets = @{
    "https://releases.jquery.com" = "https://releases.jquery.com/"
}

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.Type
    $ws.Cells.Item($r.Row, 3).Value = $r.Product
    $ws.Cells.Item($r.Row, 4).Value = $r.Url
}

# --- 2. Rebuild the hyperlinks so each one follows its row -----------------
# Hyperlinks are bound to a fixed range, not to cell content, so clear all
# of them first (a single pass skips every other item because the
# collection mutates while being walked, so repeat until it is empty) and
# then re-add one per new D-column position.
for ($pass = 0; $pass -lt 6; $pass++) {
    foreach ($h in $ws.Hyperlinks) {
        $h.Delete()
    }
}

foreach ($r in $rows) {
    $target = $r.Url
    if ($hyperlinkTargets.ContainsKey($r.Url)) {
        $target = $hyperlinkTargets[$r.Url]
    }
    $cell = $ws.Cells.Item($r.Row, 4)
    $ws.Hyperlinks.Add($cell, $target)
    $cell.Style = "Hyperlink"
}

# --- 3. Column widths: split B away from the B:C group ---------------------
# Column B (Type) now has its own explicit width; column C (Product) keeps
# its previous best-fit width unchanged.
$ws.Columns("B").ColumnWidth = 17.666666666666668

# --- 4. Selection moved from F10 to E9 -------------------------------------
$ws.Range("E9").Select()
